$wb = $excel.ActiveWorkbook

# The workbook originally had two sheets:
#   "Sheet1" - a text-formatted duplicate of the contact list
#   "Sheet2" - the "real" contact list (account number / name / mail)
# The edit removes the duplicate "Sheet1" and keeps the remaining sheet,
# renamed to "Sheet1", as the workbook's only worksheet.

$oldSheet1 = $wb.Worksheets.Item("Sheet1")
[void]$oldSheet1.Delete()

$mainSheet = $wb.Worksheets.Item("Sheet2")
$mainSheet.Name = "Sheet1"

# The header in column B changes from "담당자" (contact) to "이름" (name)
$mainSheet.Range("B1").Value = "이름"

# The active cell/selection on the sheet moves to B1
[void]$mainSheet.Range("B1").Select()
